$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 183, shifting existing rows 183:194 down to 184:195.
$ws.Rows.Item(183).Insert()

# Populate the newly inserted row 183 with a new weekly Piña price record
# (same market/product metadata as the surrounding rows; new date & price figures).
$ws.Range("A183").Value = 4
$ws.Range("B183").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C183").Value = "Los Lagos"
$ws.Range("D183").Value = 44585
$ws.Range("E183").Value = 10
$ws.Range("F183").Value = "Fruta"
$ws.Range("G183").Value = 100108
$ws.Range("H183").Value = "Tropicales y subtropicales"
$ws.Range("I183").Value = 100108005
$ws.Range("J183").Value = "Piña"
$ws.Range("K183").Value = "Caramelo"
$ws.Range("L183").Value = "Tercera"
$ws.Range("M183").Value = 120
$ws.Range("N183").Value = 19000
$ws.Range("O183").Value = 19000
$ws.Range("P183").Value = 19000
$ws.Range("Q183").Value = "$/caja 16 unidades"
$ws.Range("R183").Value = "Ecuador"
$ws.Range("S183").Value = 1188
$ws.Range("T183").Value = 16
